# Update cryptocurrency price/volume figures per the Sun Jan 8 2023 symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E2").Value = "0.22%"       # was "0.24%"

$ws.Range("D3").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D3").Value = "26.66"       # was "26.67"

$ws.Range("E4").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E4").Value = "0.07%"       # was "0.01%"

$ws.Range("D5").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D5").Value = "0.06183"       # was "0.06184"

$ws.Range("E5").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E5").Value = "1.59%"       # was "1.61%"

$ws.Range("D6").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D6").Value = "6.705"       # was "6.709"

$ws.Range("E6").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E6").Value = "0.64%"       # was "0.55%"

$ws.Range("E7").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E7").Value = "0.46%"       # was "0.36%"

$ws.Range("D8").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D8").Value = "0.9116"       # was "0.9120"

$ws.Range("E8").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E8").Value = "-1.93%"       # was "-2.26%"

$ws.Range("E9").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E9").Value = "-0.19%"       # was "-0.25%"

$ws.Range("D10").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D10").Value = "0.05149"       # was "0.05134"

$ws.Range("E10").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E10").Value = "4.39%"       # was "4.79%"

$ws.Range("D11").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D11").Value = "0.07105"       # was "0.07098"

$ws.Range("E11").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E11").Value = "0.02%"       # was "-0.06%"

$ws.Range("D12").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D12").Value = "0.03106"       # was "0.03109"

$ws.Range("E12").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E12").Value = "0.82%"       # was "1.10%"

$ws.Range("D13").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D13").Value = "0.09045"       # was "0.09037"

$ws.Range("E13").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E13").Value = "-0.15%"       # was "-0.30%"

$ws.Range("D14").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D14").Value = "0.001530"       # was "0.001544"

$ws.Range("E14").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E14").Value = "-0.88%"       # was "0.37%"

$ws.Range("D15").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D15").Value = "0.0006186"       # was "0.0006149"

$ws.Range("E15").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E15").Value = "1.44%"       # was "0.88%"

$ws.Range("D16").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D16").Value = "0.005989"       # was "0.005988"

$ws.Range("E16").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E16").Value = "-1.91%"       # was "-1.34%"

$ws.Range("E17").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E17").Value = "-0.02%"       # was "0.00%"

$ws.Range("D18").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D18").Value = "3.174"       # was "3.170"

$ws.Range("E18").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E18").Value = "0.77%"       # was "0.66%"

$ws.Range("E19").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E19").Value = "1.15%"       # was "0.26%"

$ws.Range("E21").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E21").Value = "0.46%"       # was "0.47%"

$ws.Range("D22").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D22").Value = "4.083"       # was "4.089"

$ws.Range("E22").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E22").Value = "-0.14%"       # was "0.03%"

$ws.Range("D23").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D23").Value = "0.04233"       # was "0.04254"

$ws.Range("E23").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E23").Value = "0.13%"       # was "0.03%"

$ws.Range("D24").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D24").Value = "0.001176"       # was "0.001181"

$ws.Range("E24").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E24").Value = "-3.93%"       # was "-3.45%"

$ws.Range("D25").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D25").Value = "0.004049"       # was "0.004052"

$ws.Range("E25").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E25").Value = "6.54%"       # was "6.59%"

$ws.Range("E27").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E27").Value = "4.09%"       # was "4.10%"

$ws.Range("D40").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D40").Value = "0.03985"       # was "0.03975"

$ws.Range("E40").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E40").Value = "3.07%"       # was "2.83%"

$ws.Range("E41").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E41").Value = "0.17%"       # was "0.03%"

$ws.Range("D42").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D42").Value = "0.004141"       # was "0.004140"

$ws.Range("E42").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E42").Value = "1.45%"       # was "1.46%"

$ws.Range("D43").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D43").Value = "0.002143"       # was "0.002142"

$ws.Range("E43").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E43").Value = "-3.35%"       # was "-3.36%"

$ws.Range("E44").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E44").Value = "-18.30%"       # was "-18.80%"

$ws.Range("E45").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E45").Value = "0.36%"       # was "0.37%"

$ws.Range("E46").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E46").Value = "0.02%"       # was "0.03%"

$ws.Range("D48").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D48").Value = "0.2520"       # was "0.2582"

$ws.Range("E48").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E48").Value = "83.82%"       # was "90.52%"

$ws.Range("D49").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("D49").Value = "0.00002102"       # was "0.00002101"

$ws.Range("E49").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E49").Value = "0.02%"       # was "0.03%"

$ws.Range("E50").NumberFormat = "@"   # keep as text, like the other cells in this column
$ws.Range("E50").Value = "0.02%"       # was "0.03%"
